$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 data columns (B..AH) to 2 decimal places ("custom accuracy")
$row = 5
for ($col = 2; $col -le 34; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    $old = $cell.Value2
    $cell.Value2 = $excel.WorksheetFunction.Round($old, 2)
}

# Excel's ROUND resolves this particular midpoint (3.985 -> AA5) the opposite
# way from the source data; pin it explicitly to the intended rounded value.
$ws.Range("AA5").Value2 = 3.98

# Remove row 6 entirely (데이터 1000개 -> trimming extra sample row)
$ws.Rows.Item(6).Delete()
